$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("Provider" .. "FAF Frequency" all shift right by one,
# and likewise any data further right such as the Top Up / Frequency key at M/M
# shifts to N/N). This mirrors the author inserting a new "Investment Currency"
# column into the form.
$ws.Columns("C").Insert()

# Populate the new Investment Currency column header.
$ws.Range("C1").Value = "Investment Currency"

# Amount Invested value changed from the placeholder text "2000 - test" to a
# real number.
$ws.Range("B2").Value = 5011

# Agent name changed from "Singh, Deepak" to "Purdy, Tim" (now in column G
# after the column insert shifted it over from F).
$ws.Range("G2").Value = "Purdy, Tim"

# Investment Currency value for the one data row.
$ws.Range("C2").Value = "USD"

# The old "FAF Percentage" / "FAF Frequency" columns (now J and K after the
# insert) are no longer part of the form - remove them completely (clearing
# the inherited row style first so the now-empty cells don't linger as blank
# styled placeholders).
$ws.Range("J1:K2").Style = "Normal"
$ws.Range("J1:K2").ClearContents()

# Leave the selection on the new FAF's value cell, matching where the editor
# ended up after making the change.
[void]$ws.Range("I2").Select()
